# Updated cryptos list (Price / Volume(1h) refresh + a few row reorderings)
# For D-column values that look numeric, force text formatting first so the
# cell keeps its original "text string" representation (e.g. "570.52",
# "0.999", "10.00") instead of Excel auto-coercing it to a number/double
# (which would also silently drop meaningful trailing zeros). The style is
# reset back to "Normal" afterward so no stray NumberFormat/style lingers
# on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.604.72"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "2.944.41"
$ws.Range("E3").Value = "  -2.15%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").Value = "2.942.18"
$ws.Range("E9").Value = "  -2.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.61%  "
$ws.Range("E11").Value = "  -3.70%  "
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("E13").Value = "  -2.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").Value = "65.603.95"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").Value = "3.435.27"
$ws.Range("E17").Value = "  -2.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("D19").Value = "2.942.53"
$ws.Range("E19").Value = "  -2.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +12.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "446.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.695"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.02%  "
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("E32").Value = "  -3.64%  "
$ws.Range("E33").Value = "  +4.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.21"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E36").Value = "  -1.79%  "
$ws.Range("E37").Value = "  -2.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "46.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.15%  "
$ws.Range("E40").Value = "  -9.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.302"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.120"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "384.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.88%  "
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("D47").Value = "2.673.62"
$ws.Range("E47").Value = "  -4.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E51").Value = "  +1.55%  "
